$d = $word.ActiveDocument

$replacements = @(
    @{ old = "Distance entre le domicile et le lieu de soins"; new = "Distance between home and place of care" },
    @{ old = "Problème de transport"; new = "Transport problem" },
    @{ old = "Manque de personnel compétent pour vous administrer les soins dont vous avez besoin"; new = "Lack of qualified staff to provide the care you need" },
    @{ old = "Délais d'attente très longs"; new = "Very long waiting times" },
    @{ old = "Difficultés pour payer les soins"; new = "Difficulties paying for care" },
    @{ old = "Peur des examens médicaux, de l'hôpital, ou d'autre chose"; new = "Fear of medical tests, the hospital, or something else" },
    @{ old = "Manque de temps à cause du travail, des enfants, ou autre"; new = "Lack of time because of work, children, or other" },
    @{ old = "Manque d'information"; new = "Lack of information" },
    @{ old = "Difficultés liées à la langue"; new = "Language difficulties" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
